$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new line rows (line7, line8) were inserted into the data, the
# former extr1..extr8 rows shifted down two rows (becoming extr1..extr6 in
# rows 10-15) and two brand-new extr rows (extr7, extr8) were appended at
# the bottom (rows 16-17). Several C/D/E values were also updated.

# Make sure rows 16 and 17 exist with the same formatting (bold, bordered,
# centered) as the rest of column A before writing values into them.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A16:A17").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$data = @(
  @{Row=8;  A=6;  B="line7"; C=14; D=11; E=$true}
  @{Row=9;  A=7;  B="line8"; C=16; D=9;  E=$true}
  @{Row=10; A=8;  B="extr1"; C=5;  D=12; E=$true}
  @{Row=11; A=9;  B="extr2"; C=5;  D=9;  E=$true}
  @{Row=12; A=10; B="extr3"; C=10; D=11; E=$false}
  @{Row=13; A=11; B="extr4"; C=7;  D=8;  E=$true}
  @{Row=14; A=12; B="extr5"; C=9;  D=11; E=$false}
  @{Row=15; A=13; B="extr6"; C=7;  D=11; E=$false}
  @{Row=16; A=14; B="extr7"; C=5;  D=7;  E=$true}
  @{Row=17; A=15; B="extr8"; C=8;  D=5;  E=$false}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
}
